$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (un_franzosa_ControlvsCD_Fp) for the
# un_franzosa_ControlvsCD_ConvCD data, shifting subsequent rows down.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.4
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6
$ws.Range("H9").Value = 0.6

# Insert a new row above row 15 (now un_franzosa_ControlvsUC_Fp after the
# previous insert) for the un_franzosa_ControlvsUC_ConvUC data, shifting
# subsequent rows down.
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8
$ws.Range("H15").Value = 0.8
